# Add this week's meeting attendance row (row 9) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date for the new meeting (01/11/23 -> serial 45231)
$ws.Range("A9").Value = 45231
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B9").Value = "Yes"
$ws.Range("C9").Value = "Yes"
$ws.Range("D9").Value = "Yes"
$ws.Range("E9").Value = "Yes"

# Update the selection to match the recorded state in the diff
$ws.Range("F9").Select()
